# Financial statement figure updates (yearly revision pass).
# Applies the refreshed Income Statement / Balance Sheet / Cash Flow
# figures for the RDY sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8: Total Revenue
$ws.Range("D8").Value = 2053700
$ws.Range("E8").Value = 2036100
$ws.Range("F8").Value = 2237100
$ws.Range("G8").Value = 2142800
$ws.Range("H8").Value = 1911200
$ws.Range("I8").Value = 1681200
$ws.Range("J8").Value = 1398800

# Row 9: Cost of Revenue
$ws.Range("D9").Value = 950400
$ws.Range("E9").Value = 903100
$ws.Range("F9").Value = 902700
$ws.Range("G9").Value = 907900
$ws.Range("H9").Value = 815100
$ws.Range("I9").Value = 805200
$ws.Range("J9").Value = 1256100

# Row 10: Gross Profit
$ws.Range("D10").Value = 1103400
$ws.Range("E10").Value = 1133000
$ws.Range("F10").Value = 1334400
$ws.Range("G10").Value = 1234900
$ws.Range("H10").Value = 1096100
$ws.Range("I10").Value = 876000
$ws.Range("J10").Value = 142800

# Row 12: Research Development
$ws.Range("D12").Value = 246300
$ws.Range("E12").Value = 264400
$ws.Range("F12").Value = 240400
$ws.Range("G12").Value = 249800
$ws.Range("H12").Value = 179300
$ws.Range("I12").Value = 111000
$ws.Range("J12").Value = 85500

# Row 14: Non Recurring
$ws.Range("D14").Value = 800
$ws.Range("E14").Value = 1600
$ws.Range("F14").Value = 2800
$ws.Range("G14").Value = 9100
$ws.Range("H14").Value = -7200
$ws.Range("I14").Value = 7500
$ws.Range("J14").Value = 15200

# Row 15: Others
$ws.Range("D15").Value = 74000
$ws.Range("E15").Value = 74600
$ws.Range("F15").Value = 72400
$ws.Range("G15").Value = 34400
$ws.Range("J15").Value = 22900

# Row 17: Total Operating Expenses
$ws.Range("D17").Value = 1881400
$ws.Range("E17").Value = 1840900
$ws.Range("F17").Value = 1808800
$ws.Range("G17").Value = 1762700
$ws.Range("H17").Value = 1534800
$ws.Range("I17").Value = 1375900
$ws.Range("J17").Value = 1134900

# Row 18: Operating Income or Loss
$ws.Range("D18").Value = 172300
$ws.Range("E18").Value = 195200
$ws.Range("F18").Value = 428300
$ws.Range("G18").Value = 380100
$ws.Range("H18").Value = 376400
$ws.Range("I18").Value = 305300
$ws.Range("J18").Value = 263900

# Row 20: Total Other Income/Expenses Net
$ws.Range("D20").Value = 46400
$ws.Range("E20").Value = 25900
$ws.Range("F20").Value = -23900
$ws.Range("G20").Value = 42900
$ws.Range("H20").Value = 26700
$ws.Range("I20").Value = 22900
$ws.Range("J20").Value = 18500

# Row 21: Earnings Before Interest And Taxes
$ws.Range("D21").Value = 388200
$ws.Range("E21").Value = 384200
$ws.Range("F21").Value = 552700
$ws.Range("G21").Value = 540200
$ws.Range("H21").Value = 506000
$ws.Range("I21").Value = 408400
$ws.Range("J21").Value = 357900

# Row 22: Interest Expense
$ws.Range("D22").Value = 11400
$ws.Range("E22").Value = 9200
$ws.Range("F22").Value = 11900
$ws.Range("G22").Value = 15800
$ws.Range("H22").Value = 18400
$ws.Range("I22").Value = 14700
$ws.Range("J22").Value = 15400

# Row 23: Income Before Tax
$ws.Range("D23").Value = 207400
$ws.Range("E23").Value = 211900
$ws.Range("F23").Value = 392400
$ws.Range("G23").Value = 407200
$ws.Range("H23").Value = 384700
$ws.Range("I23").Value = 313400
$ws.Range("J23").Value = 267000

# Row 24: Income Tax Expense
$ws.Range("D24").Value = 46700
$ws.Range("E24").Value = 37800
$ws.Range("F24").Value = 103100
$ws.Range("G24").Value = 86500
$ws.Range("H24").Value = 73700
$ws.Range("I24").Value = 70900
$ws.Range("J24").Value = 60800

# Row 26: Income After Tax
$ws.Range("D26").Value = 160700
$ws.Range("E26").Value = 174100
$ws.Range("F26").Value = 289400
$ws.Range("G26").Value = 320700
$ws.Range("H26").Value = 311100
$ws.Range("I26").Value = 242600
$ws.Range("J26").Value = 206200

# Row 27: Net Income From Continuing Ops
$ws.Range("D27").Value = 160700
$ws.Range("E27").Value = 174100
$ws.Range("F27").Value = 289400
$ws.Range("G27").Value = 320700
$ws.Range("H27").Value = 311100
$ws.Range("I27").Value = 242600
$ws.Range("J27").Value = 206200

# Row 29: Discontinued Operations
$ws.Range("D29").Value = -18900

# Row 32: Other Items
$ws.Range("D32").Value = -46400
$ws.Range("E32").Value = -25900
$ws.Range("F32").Value = 23900
$ws.Range("G32").Value = -42900
$ws.Range("H32").Value = -26700
$ws.Range("I32").Value = -22900
$ws.Range("J32").Value = -18500

# Row 33: Net Income
$ws.Range("D33").Value = 141800
$ws.Range("E33").Value = 174100
$ws.Range("F33").Value = 289400
$ws.Range("G33").Value = 320700
$ws.Range("H33").Value = 311100
$ws.Range("I33").Value = 242600
$ws.Range("J33").Value = 206200

# Row 35: Net Income Applicable To Common Shares
$ws.Range("D35").Value = 141800
$ws.Range("E35").Value = 174100
$ws.Range("F35").Value = 289400
$ws.Range("G35").Value = 320700
$ws.Range("H35").Value = 311100
$ws.Range("I35").Value = 242600
$ws.Range("J35").Value = 206200

# Row 41: Cash And Cash Equivalents
$ws.Range("D41").Value = 35700
$ws.Range("E41").Value = 53300
$ws.Range("F41").Value = 71200
$ws.Range("G41").Value = 78000
$ws.Range("H41").Value = 122200
$ws.Range("I41").Value = 74300
$ws.Range("J41").Value = 66400

# Row 42: Short Term Investments
$ws.Range("D42").Value = 265100
$ws.Range("E42").Value = 206300
$ws.Range("F42").Value = 506600
$ws.Range("G42").Value = 495400
$ws.Range("H42").Value = 362700
$ws.Range("I42").Value = 246700
$ws.Range("J42").Value = 193400

# Row 43: Net Receivables
$ws.Range("D43").Value = 628400
$ws.Range("E43").Value = 653200
$ws.Range("F43").Value = 688000
$ws.Range("G43").Value = 666200
$ws.Range("H43").Value = 477700
$ws.Range("I43").Value = 924600
$ws.Range("J43").Value = 366400

# Row 44: Inventory
$ws.Range("D44").Value = 420600
$ws.Range("E44").Value = 412500
$ws.Range("F44").Value = 369900
$ws.Range("G44").Value = 369100
$ws.Range("H44").Value = 346900
$ws.Range("I44").Value = 312300
$ws.Range("J44").Value = 279800

# Row 45: Other Current Assets
$ws.Range("D45").Value = 235700
$ws.Range("E45").Value = 126000
$ws.Range("F45").Value = 95100
$ws.Range("G45").Value = 124100
$ws.Range("H45").Value = 190600
$ws.Range("I45").Value = 296100
$ws.Range("J45").Value = 105400

# Row 46: Total Current Assets
$ws.Range("D46").Value = 1585500
$ws.Range("E46").Value = 1451400
$ws.Range("F46").Value = 1730700
$ws.Range("G46").Value = 1732900
$ws.Range("H46").Value = 1500200
$ws.Range("I46").Value = 1239400
$ws.Range("J46").Value = 1011500

# Row 47: Long Term Investments
$ws.Range("D47").Value = 69700
$ws.Range("E47").Value = 101900
$ws.Range("F47").Value = 47700
$ws.Range("G47").Value = 55700
$ws.Range("H47").Value = 11700
$ws.Range("I47").Value = 9800
$ws.Range("J47").Value = 10600

# Row 48: Property Plant and Equipment
$ws.Range("D48").Value = 836800
$ws.Range("E48").Value = 826500
$ws.Range("F48").Value = 780300
$ws.Range("G48").Value = 695400
$ws.Range("H48").Value = 642400
$ws.Range("I48").Value = 1093600
$ws.Range("J48").Value = 480700

# Row 49: Goodwill
$ws.Range("D49").Value = 702900
$ws.Range("E49").Value = 703900
$ws.Range("F49").Value = 356400
$ws.Range("G49").Value = 237600
$ws.Range("H49").Value = 212500
$ws.Range("I49").Value = 202700
$ws.Range("J49").Value = 195600

# Row 52: Other Assets
$ws.Range("D52").Value = 67400
$ws.Range("E52").Value = 94900
$ws.Range("F52").Value = 87600
$ws.Range("G52").Value = 94800
$ws.Range("H52").Value = 94700
$ws.Range("I52").Value = 59800
$ws.Range("J52").Value = 29200

# Row 54: Total Assets
$ws.Range("D54").Value = 3262200
$ws.Range("E54").Value = 3178600
$ws.Range("F54").Value = 3002600
$ws.Range("G54").Value = 2816300
$ws.Range("H54").Value = 2461400
$ws.Range("I54").Value = 2058700
$ws.Range("J54").Value = 1727600

# Row 57: Accounts Payable
$ws.Range("E57").Value = 194000
$ws.Range("F57").Value = 177900
$ws.Range("G57").Value = 154100
$ws.Range("H57").Value = 151900
$ws.Range("I57").Value = 171500
$ws.Range("J57").Value = 137400

# Row 58: Short/Current Long Term Debt
$ws.Range("D58").Value = 370500
$ws.Range("E58").Value = 632400
$ws.Range("F58").Value = 330100
$ws.Range("G58").Value = 416700
$ws.Range("H58").Value = 347100
$ws.Range("I58").Value = 349000
$ws.Range("J58").Value = 229600

# Row 59: Other Current Liabilities
$ws.Range("D59").Value = 637000
$ws.Range("E59").Value = 402700
$ws.Range("F59").Value = 426800
$ws.Range("G59").Value = 354500
$ws.Range("H59").Value = 282800
$ws.Range("I59").Value = 261600
$ws.Range("J59").Value = 261500

# Row 60: Total Current Liabilities
$ws.Range("D60").Value = 1007700
$ws.Range("E60").Value = 1229100
$ws.Range("F60").Value = 934800
$ws.Range("G60").Value = 925400
$ws.Range("H60").Value = 781800
$ws.Range("I60").Value = 782200
$ws.Range("J60").Value = 628400

# Row 61: Long Term Debt
$ws.Range("D61").Value = 362800
$ws.Range("E61").Value = 78800
$ws.Range("F61").Value = 154500
$ws.Range("G61").Value = 206900
$ws.Range("H61").Value = 299900
$ws.Range("I61").Value = 182600
$ws.Range("J61").Value = 236200

# Row 62: Other Liabilities
$ws.Range("D62").Value = 63100
$ws.Range("E62").Value = 77000
$ws.Range("F62").Value = 57600
$ws.Range("G62").Value = 74600
$ws.Range("H62").Value = 66800
$ws.Range("I62").Value = 44200
$ws.Range("J62").Value = 32400

# Row 66: Total Liabilities
$ws.Range("D66").Value = 1433600
$ws.Range("E66").Value = 1384900
$ws.Range("F66").Value = 1146900
$ws.Range("G66").Value = 1206800
$ws.Range("H66").Value = 1148400
$ws.Range("I66").Value = 1006200
$ws.Range("J66").Value = 897000

# Row 72: Retained Earnings
$ws.Range("D72").Value = 1663800
$ws.Range("E72").Value = 1579400
$ws.Range("F72").Value = 1455400
$ws.Range("G72").Value = 1225100
$ws.Range("H72").Value = 955100
$ws.Range("I72").Value = 685900
$ws.Range("J72").Value = 481400

# Row 76: Total Stockholder Equity
$ws.Range("D76").Value = 1828600
$ws.Range("E76").Value = 1793700
$ws.Range("F76").Value = 1855700
$ws.Range("G76").Value = 1609400
$ws.Range("H76").Value = 1313000
$ws.Range("I76").Value = 1052500
$ws.Range("J76").Value = 830600

# Row 81: Net Income
$ws.Range("D81").Value = 141800
$ws.Range("E81").Value = 174100
$ws.Range("F81").Value = 289400
$ws.Range("G81").Value = 320700
$ws.Range("H81").Value = 311100
$ws.Range("I81").Value = 242600
$ws.Range("J81").Value = 206200

# Row 83: Depreciation
$ws.Range("D83").Value = 169300
$ws.Range("E83").Value = 163100
$ws.Range("F83").Value = 148200
$ws.Range("G83").Value = 117100
$ws.Range("H83").Value = 102800
$ws.Range("I83").Value = 80200
$ws.Range("J83").Value = 75400

# Row 89: Total Cash Flow From Operating Activities
$ws.Range("D89").Value = 260700
$ws.Range("E89").Value = 311100
$ws.Range("F89").Value = 596400
$ws.Range("G89").Value = 362000
$ws.Range("H89").Value = 281400
$ws.Range("I89").Value = 192600
$ws.Range("J89").Value = 233500

# Row 91: Capital Expenditures
$ws.Range("D91").Value = -134300
$ws.Range("E91").Value = -177500
$ws.Range("F91").Value = -173800
$ws.Range("G91").Value = -135000
$ws.Range("H91").Value = -145800
$ws.Range("I91").Value = -96400
$ws.Range("J91").Value = -99200

# Row 94: Total Cash Flows From Investing Activities
$ws.Range("D94").Value = -215200
$ws.Range("E94").Value = -267100
$ws.Range("F94").Value = -295300
$ws.Range("G94").Value = -331200
$ws.Range("H94").Value = -240300
$ws.Range("I94").Value = -201600
$ws.Range("J94").Value = -269900

# Row 96: Dividends Paid
$ws.Range("D96").Value = -57700
$ws.Range("E96").Value = -49000
$ws.Range("F96").Value = -59400
$ws.Range("G96").Value = -51900
$ws.Range("H96").Value = -43200
$ws.Range("I96").Value = -39200
$ws.Range("J96").Value = -32000

# Row 100: Total Cash Flows From Financing Activities
$ws.Range("D100").Value = -64200
$ws.Range("E100").Value = -53400
$ws.Range("F100").Value = -245800
$ws.Range("G100").Value = -59500
$ws.Range("H100").Value = -3100
$ws.Range("I100").Value = -25900
$ws.Range("J100").Value = 54000

# Row 101: Effect Of Exchange Rate Changes 
$ws.Range("E101").Value = -7100
$ws.Range("F101").Value = -62100
$ws.Range("G101").Value = -15400
$ws.Range("H101").Value = 11100
$ws.Range("I101").Value = 1400
$ws.Range("J101").Value = 7200

# Row 102: Change In Cash and Cash Equivalents 
$ws.Range("D102").Value = -17900
$ws.Range("E102").Value = -16500
$ws.Range("F102").Value = -6800
$ws.Range("G102").Value = -44200
$ws.Range("H102").Value = 49100
$ws.Range("I102").Value = -33600
$ws.Range("J102").Value = 24900
